$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1658
$ws1.Range("F9").Value = 599
$ws1.Range("F14").Value = 233

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6295
$ws3.Range("F5").Value = 46

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6295
$ws4.Range("F6").Value = 46
$ws4.Range("F16").Value = 1658
$ws4.Range("F24").Value = 599
$ws4.Range("F37").Value = 233
